$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-21 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-22 Monday", 2) | Out-Null
$d.Content.Find.Execute("985×8=7880", $true, $false, $false, $false, $false, $true, 1, $false, "908×4=3632", 2) | Out-Null
$d.Content.Find.Execute("849×8=6792", $true, $false, $false, $false, $false, $true, 1, $false, "856×5=4280", 2) | Out-Null
$d.Content.Find.Execute("581×6=3486", $true, $false, $false, $false, $false, $true, 1, $false, "956×3=2868", 2) | Out-Null
$d.Content.Find.Execute("309×7=2163", $true, $false, $false, $false, $false, $true, 1, $false, "785×9=7065", 2) | Out-Null
$d.Content.Find.Execute("830×7=5810", $true, $false, $false, $false, $false, $true, 1, $false, "461×6=2766", 2) | Out-Null
$d.Content.Find.Execute("224×2=448", $true, $false, $false, $false, $false, $true, 1, $false, "293×2=586", 2) | Out-Null
$d.Content.Find.Execute("574×3=1722", $true, $false, $false, $false, $false, $true, 1, $false, "678×8=5424", 2) | Out-Null
$d.Content.Find.Execute("863×2=1726", $true, $false, $false, $false, $false, $true, 1, $false, "332×9=2988", 2) | Out-Null
$d.Content.Find.Execute("371×9=3339", $true, $false, $false, $false, $false, $true, 1, $false, "429×4=1716", 2) | Out-Null
$d.Content.Find.Execute("183×2=366", $true, $false, $false, $false, $false, $true, 1, $false, "851×9=7659", 2) | Out-Null
$d.Content.Find.Execute("139×3=417", $true, $false, $false, $false, $false, $true, 1, $false, "995×3=2985", 2) | Out-Null
$d.Content.Find.Execute("554×3=1662", $true, $false, $false, $false, $false, $true, 1, $false, "668×4=2672", 2) | Out-Null
$d.Content.Find.Execute("985×2=1970", $true, $false, $false, $false, $false, $true, 1, $false, "859×4=3436", 2) | Out-Null
$d.Content.Find.Execute("520×3=1560", $true, $false, $false, $false, $false, $true, 1, $false, "540×9=4860", 2) | Out-Null
$d.Content.Find.Execute("778×7=5446", $true, $false, $false, $false, $false, $true, 1, $false, "660×6=3960", 2) | Out-Null
$d.Content.Find.Execute("824×3=2472", $true, $false, $false, $false, $false, $true, 1, $false, "750×5=3750", 2) | Out-Null
$d.Content.Find.Execute("981×3=2943", $true, $false, $false, $false, $false, $true, 1, $false, "358×4=1432", 2) | Out-Null
$d.Content.Find.Execute("951×2=1902", $true, $false, $false, $false, $false, $true, 1, $false, "689×6=4134", 2) | Out-Null
$d.Content.Find.Execute("612×7=4284", $true, $false, $false, $false, $false, $true, 1, $false, "221×4=884", 2) | Out-Null
$d.Content.Find.Execute("200×7=1400", $true, $false, $false, $false, $false, $true, 1, $false, "779×7=5453", 2) | Out-Null
$d.Content.Find.Execute("622×7=4354", $true, $false, $false, $false, $false, $true, 1, $false, "211×7=1477", 2) | Out-Null
$d.Content.Find.Execute("375×6=2250", $true, $false, $false, $false, $false, $true, 1, $false, "705×5=3525", 2) | Out-Null
$d.Content.Find.Execute("436×4=1744", $true, $false, $false, $false, $false, $true, 1, $false, "966×3=2898", 2) | Out-Null
$d.Content.Find.Execute("188×9=1692", $true, $false, $false, $false, $false, $true, 1, $false, "605×2=1210", 2) | Out-Null
$d.Content.Find.Execute("631×5=3155", $true, $false, $false, $false, $false, $true, 1, $false, "234×5=1170", 2) | Out-Null
